$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before assigning so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00") into numbers, then restore
# the original (default) style afterwards so no stray cell styles are introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "41.929.79"
$ws.Range("E2").Value = "  -1.96%  "

# Row 3
$ws.Range("D3").Value = "2.241.58"
$ws.Range("E3").Value = "  -2.25%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "246.96"
$ws.Range("E5").Value = "  -1.95%  "

# Row 6
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -2.71%  "

# Row 7
$ws.Range("D7").Value = "74.84"
$ws.Range("E7").Value = "  +0.83%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  -2.70%  "

# Row 10
$ws.Range("D10").Value = "40.19"
$ws.Range("E10").Value = "  +1.43%  "

# Row 11
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  -3.62%  "

# Row 12
$ws.Range("D12").Value = "7.13"
$ws.Range("E12").Value = "  -3.06%  "

# Row 13
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
$ws.Range("D14").Value = "2.583.16"
$ws.Range("E14").Value = "  -2.09%  "

# Row 15
$ws.Range("D15").Value = "14.80"
$ws.Range("E15").Value = "  -3.67%  "

# Row 16
$ws.Range("D16").Value = "0.860"
$ws.Range("E16").Value = "  -2.00%  "

# Row 17
$ws.Range("D17").Value = "2.238.10"
$ws.Range("E17").Value = "  -1.97%  "

# Row 18
$ws.Range("D18").Value = "41.915.66"
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0980"
$ws.Range("E19").Value = "  -2.21%  "

# Row 20
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  -2.70%  "

# Row 21
$ws.Range("D21").Value = "71.64"
$ws.Range("E21").Value = "  -1.27%  "

# Row 22
$ws.Range("E22").Value = "  +0.77%  "

# Row 23
$ws.Range("D23").Value = "229.66"
$ws.Range("E23").Value = "  -1.51%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.48"
$ws.Range("E24").Value = "  -0.87%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.04%  "

# Row 26
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -5.78%  "

# Row 27
$ws.Range("E27").Value = "  -4.72%  "

# Row 28
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  +11.53%  "

# Row 29
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -1.40%  "

# Row 30
$ws.Range("D30").Value = "169.04"
$ws.Range("E30").Value = "  +1.05%  "

# Row 31
$ws.Range("D31").Value = "20.54"
$ws.Range("E31").Value = "  -2.73%  "

# Row 32
$ws.Range("D32").Value = "33.81"
$ws.Range("E32").Value = "  +5.04%  "

# Row 33
$ws.Range("D33").Value = "0.0846"
$ws.Range("E33").Value = "  +3.02%  "

# Row 34
$ws.Range("D34").Value = "0.121"
$ws.Range("E34").Value = "  -4.99%  "

# Row 35
$ws.Range("D35").Value = "0.127"
$ws.Range("E35").Value = "  -0.22%  "

# Row 36
$ws.Range("D36").Value = "4.57"
$ws.Range("E36").Value = "  -3.60%  "

# Row 37
$ws.Range("E37").Value = "  +2.42%  "

# Row 38
$ws.Range("D38").Value = "0.0300"
$ws.Range("E38").Value = "  -2.29%  "

# Row 39
$ws.Range("D39").Value = "13.47"
$ws.Range("E39").Value = "  -5.69%  "

# Row 40
$ws.Range("D40").Value = "5.93"
$ws.Range("E40").Value = "  -0.88%  "

# Row 41
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").Value = "  -6.80%  "

# Row 42
$ws.Range("D42").Value = "112.17"
$ws.Range("E42").Value = "  +14.13%  "

# Row 43
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").Value = "  -4.60%  "

# Row 44
$ws.Range("D44").Value = "60.50"
$ws.Range("E44").Value = "  -2.64%  "

# Row 45
$ws.Range("D45").Value = "8.80"
$ws.Range("E45").Value = "  -3.65%  "

# Row 46
$ws.Range("E46").Value = "  -2.38%  "

# Row 47
$ws.Range("E47").Value = "  -0.67%  "

# Row 48
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").Value = "  -4.01%  "

# Row 49
$ws.Range("D49").Value = "1.18"
$ws.Range("E49").Value = "  -1.58%  "

# Row 50
$ws.Range("D50").Value = "4.32"
$ws.Range("E50").Value = "  -10.98%  "

# Row 51
$ws.Range("D51").Value = "4.21"
$ws.Range("E51").Value = "  -1.51%  "

# Restore default style for the price column so only values changed, not formats.
$priceRange.Style = "Normal"
